$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F2").Value = "26 jun. 2023, 18:25:57"
